$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseXpath = "/Envelope/Body/GetInfoByStateResponse/GetInfoByStateResult/NewDataSet/Table[2]/"

$fields = @("CITY", "STATE", "ZIP", "AREA_CODE", "TIME_ZONE")
$values = @("[A-Z a-z].*", "[A-Z]{2}", "[0-9]{5}", "[0-9]{3}", "[A-Z]{1}")

$row = 7
for ($i = 0; $i -lt $fields.Length; $i++) {
    $ws.Cells.Item($row, 2).Value = $baseXpath + $fields[$i]
    $ws.Cells.Item($row, 3).Value = $values[$i]
    $row++
}

$ws.Range("B12").Select() | Out-Null
